$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C holds "Fitness" values keyed by row = Generation (col B) + 2.
# Apply the corrected fitness values for the affected generation ranges.

for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 3).Value = 7723
}

for ($r = 7; $r -le 41; $r++) {
    $ws.Cells.Item($r, 3).Value = 7534
}

for ($r = 42; $r -le 43; $r++) {
    $ws.Cells.Item($r, 3).Value = 7320
}

for ($r = 44; $r -le 50; $r++) {
    $ws.Cells.Item($r, 3).Value = 7318
}

for ($r = 141; $r -le 252; $r++) {
    $ws.Cells.Item($r, 3).Value = 7293
}
